# cfb_weather.xlsx update - refresh odds line-movement columns (Open/Current/Move_s),
# a handful of corrected wind-direction readings, and the run Timestamp.

$wb      = $excel.ActiveWorkbook
$wsFBS   = $wb.Worksheets.Item("FBS")
$wsOther = $wb.Worksheets.Item("Other")

# ---------------------------------------------------------------------------
# Corrected wind_dir_fg / gs_fg readings (forecast-game wind direction)
# ---------------------------------------------------------------------------
$wsFBS.Range("Q4").Value   = "ENE"
$wsFBS.Range("Q19").Value  = "NW"
$wsFBS.Range("Q35").Value  = "WSW"

$wsOther.Range("S9").Value  = "NNW"
$wsOther.Range("S33").Value = "SW"
$wsOther.Range("S41").Value = "N"

# ---------------------------------------------------------------------------
# Freshly pulled opening/current spread (Open / Current) for each FBS game,
# with the spread-movement tracker (Move_s) reset to 0 for the new pull.
# ---------------------------------------------------------------------------
$openLines = @{
    2  = 1.5
    3  = 14.5
    4  = 2.5
    5  = 3.5
    6  = -13.5
    7  = -5.5
    8  = -7.5
    9  = -5.5
    10 = -2.5
    11 = -7
    12 = -7.5
    13 = -19.5
    14 = -10.5
    15 = 2
    16 = 6.5
    17 = 8.5
    18 = 1
    19 = 6.5
    20 = -10.5
    21 = -11
    22 = -17.5
    23 = -17.5
    24 = -8.5
    26 = 10.5
    27 = 3
    28 = 14.5
    29 = -32
    30 = 0
    31 = -3
    32 = 14.5
    33 = 2.5
    34 = -22
    35 = -15
    36 = -9
    37 = 3.5
    38 = 6.5
    39 = 3
    40 = 9
    41 = -13
    42 = 3.5
    43 = 2.5
    44 = 14
    45 = -2.5
    46 = -7.5
    47 = 3.5
    48 = -9
    49 = -9
    50 = -16.5
    51 = 7
    52 = -2
}

foreach ($row in $openLines.Keys) {
    $line = $openLines[$row]
    $wsFBS.Range("AA$row").Value = $line
    $wsFBS.Range("AB$row").Value = $line
    $wsFBS.Range("AF$row").Value = 0
}

# ---------------------------------------------------------------------------
# Stamp every FBS row with the timestamp of this refresh run.
# ---------------------------------------------------------------------------
$newTimestamp = "2025-10-07T10:36:14.420713"
for ($row = 2; $row -le 52; $row++) {
    $wsFBS.Range("AK$row").Value = $newTimestamp
}
